$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.290.44"
$ws.Range("E2").Value = "  +0.43%  "

# Row 3
$ws.Range("D3").Value = "2.629.53"
$ws.Range("E3").Value = "  +0.68%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").Value = "'601.83"
$ws.Range("E5").Value = "  +1.67%  "

# Row 6
$ws.Range("D6").Value = "'153.30"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("E8").Value = "  +2.88%  "

# Row 9
$ws.Range("D9").Value = "2.629.49"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("E10").Value = "  +6.33%  "

# Row 11
$ws.Range("E11").Value = "  +0.75%  "

# Row 12
$ws.Range("E12").Value = "  +0.19%  "

# Row 13
$ws.Range("D13").Value = "'0.351"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14
$ws.Range("D14").Value = "'27.89"
$ws.Range("E14").Value = "  +1.07%  "

# Row 15
$ws.Range("D15").Value = "3.108.65"
$ws.Range("E15").Value = "  +0.78%  "

# Row 16
$ws.Range("E16").Value = "  +1.07%  "

# Row 17
$ws.Range("D17").Value = "67.374.89"
$ws.Range("E17").Value = "  +0.57%  "

# Row 18
$ws.Range("D18").Value = "2.630.38"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19
$ws.Range("D19").Value = "'11.23"
$ws.Range("E19").Value = "  -0.02%  "

# Row 20
$ws.Range("D20").Value = "'363.70"
$ws.Range("E20").Value = "  +1.25%  "

# Row 21
$ws.Range("E21").Value = "  -3.85%  "

# Row 22
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("E23").Value = "  +5.56%  "

# Row 24
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  -0.14%  "

# Row 25
$ws.Range("D25").Value = "'10.14"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26
$ws.Range("D26").Value = "'66.37"
$ws.Range("E26").Value = "  -7.15%  "

# Row 27
$ws.Range("E27").Value = "  +0.75%  "

# Row 28
$ws.Range("E28").Value = "  +0.72%  "

# Row 29
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("D30").Value = "'579.62"
$ws.Range("E30").Value = "  -6.78%  "

# Row 31
$ws.Range("E31").Value = "  -2.83%  "

# Row 32
$ws.Range("E32").Value = "  -1.06%  "

# Row 33
$ws.Range("E33").Value = "  -0.13%  "

# Row 34
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.128"
$ws.Range("E34").Value = "  -3.19%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.13%  "

# Row 36
$ws.Range("D36").Value = "'1.54"
$ws.Range("E36").Value = "  -1.49%  "

# Row 37
$ws.Range("E37").Value = "  -0.71%  "

# Row 38
$ws.Range("D38").Value = "'158.20"
$ws.Range("E38").Value = "  +3.01%  "

# Row 39
$ws.Range("D39").Value = "'19.41"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40
$ws.Range("E40").Value = "  +0.37%  "

# Row 41
$ws.Range("D41").Value = "'5.28"
$ws.Range("E41").Value = "  -3.55%  "

# Row 42
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("E43").Value = "  +1.65%  "

# Row 44
$ws.Range("D44").Value = "'41.20"
$ws.Range("E44").Value = "  -0.33%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").Value = "'16.35"
$ws.Range("E46").Value = "  -0.84%  "

# Row 47
$ws.Range("D47").Value = "'155.87"
$ws.Range("E47").Value = "  +0.40%  "

# Row 48
$ws.Range("E48").Value = "  -2.39%  "

# Row 49
$ws.Range("E49").Value = "  -0.77%  "

# Row 50
$ws.Range("D50").Value = "'20.93"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51
$ws.Range("E51").Value = "  +0.67%  "
